# Apply crypto price/volume updates to sheet1 (ActiveSheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.312.03'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '1.789.65'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''315.60'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '''0.5344'
$ws.Range("E7").Value = '  -1.37%  '
$ws.Range("D8").Value = '''0.3757'
$ws.Range("E8").Value = '  -1.95%  '
$ws.Range("D9").Value = '''0.07476'
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("E10").Value = '  -3.70%  '
$ws.Range("D11").Value = '''1.093'
$ws.Range("E11").Value = '  -2.30%  '
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '''20.43'
$ws.Range("E13").Value = '  -3.03%  '
$ws.Range("D14").Value = '''6.090'
$ws.Range("D15").Value = '''7.243'
$ws.Range("E15").Value = '  -1.11%  '
$ws.Range("D16").Value = '1.787.20'
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").Value = '''89.21'
$ws.Range("E17").Value = '  -2.47%  '
$ws.Range("D18").Value = '''0.00001056'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = '''0.06496'
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '''17.42'
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''1.001'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '''5.928'
$ws.Range("E22").Value = '  -0.84%  '
$ws.Range("D23").Value = '28.325.98'
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").Value = '''11.09'
$ws.Range("E24").Value = '  -2.29%  '
$ws.Range("D25").Value = '''2.091'
$ws.Range("E25").Value = '  -3.89%  '
$ws.Range("D26").Value = '''158.58'
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("D27").Value = '''20.25'
$ws.Range("E27").Value = '  -1.82%  '
$ws.Range("D28").Value = '1.991.46'
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("D29").Value = '''2.293'
$ws.Range("E29").Value = '  -4.99%  '
$ws.Range("D30").Value = '''121.54'
$ws.Range("E30").Value = '  -1.81%  '
$ws.Range("D31").Value = '''1.093'
$ws.Range("E31").Value = '  -4.32%  '
$ws.Range("D32").Value = '''0.1047'
$ws.Range("E32").Value = '  +3.63%  '
$ws.Range("D33").Value = '''3.659'
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").Value = '''5.514'
$ws.Range("E34").Value = '  -3.60%  '
$ws.Range("D35").Value = '''0.2248'
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").Value = '''0.06416'
$ws.Range("E36").Value = '  +2.56%  '
$ws.Range("D37").Value = '''0.02275'
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("D38").Value = '''5.009'
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("D39").Value = '''8.502'
$ws.Range("E39").Value = '  -4.25%  '
$ws.Range("D40").Value = '''0.6154'
$ws.Range("E40").Value = '  -3.34%  '
$ws.Range("D41").Value = '''1.436'
$ws.Range("E41").Value = '  +3.81%  '
$ws.Range("D42").Value = '''1.177'
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("D43").Value = '''10.96'
$ws.Range("E43").Value = '  -4.88%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '''13.32'
$ws.Range("E45").Value = '  -1.06%  '
$ws.Range("D46").Value = '''3.661'
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("E47").Value = '  -3.52%  '
$ws.Range("D48").Value = '''125.34'
$ws.Range("E48").Value = '  +1.21%  '
$ws.Range("D49").Value = '''1.198'
$ws.Range("E49").Value = '  +4.56%  '
$ws.Range("D50").Value = '''1.928'
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("D51").Value = '''0.06847'
$ws.Range("E51").Value = '  -0.72%  '
